$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.914.78'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '2.229.85'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.39'
$ws.Range("E5").Value = '  -4.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.72'
$ws.Range("E6").Value = '  -7.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").Value = '  -1.80%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  -6.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.00'
$ws.Range("E10").Value = '  -7.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  -4.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.07'
$ws.Range("E12").Value = '  -7.53%  '
$ws.Range("E13").Value = '  -3.10%  '
$ws.Range("D14").Value = '2.568.25'
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").Value = '2.271.46'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.808'
$ws.Range("E16").Value = '  -5.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.42'
$ws.Range("E17").Value = '  -5.16%  '
$ws.Range("D18").Value = '43.674.10'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("E19").Value = '  -3.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.08'
$ws.Range("E20").Value = '  -9.37%  '
$ws.Range("E21").Value = '  -6.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.13'
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.97'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.89'
$ws.Range("E24").Value = '  -7.74%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -10.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.89'
$ws.Range("E30").Value = '  -5.26%  '
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '150.44'
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0795'
$ws.Range("E33").Value = '  -6.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.61'
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.18'
$ws.Range("E35").Value = '  +4.17%  '
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("E37").Value = '  -9.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.72'
$ws.Range("E38").Value = '  -11.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.38'
$ws.Range("E39").Value = '  -8.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.76'
$ws.Range("E40").Value = '  -10.45%  '
$ws.Range("E41").Value = '  -6.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.24'
$ws.Range("E42").Value = '  -13.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '1.729.03'
$ws.Range("E44").Value = '  -3.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '83.14'
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("E46").Value = '  -6.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.05'
$ws.Range("E47").Value = '  -3.91%  '
$ws.Range("E48").Value = '  -5.40%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '14.55'
$ws.Range("E49").Value = '  +6.63%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.03'
$ws.Range("E50").Value = '  -3.81%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.36'
$ws.Range("E51").Value = '  -10.52%  '
